$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "31.070.72"
$ws.Range("E2").Value = "  +1.34%  "

# Row 3
$ws.Range("D3").Value = "1.957.36"
$ws.Range("E3").Value = "  +0.61%  "

# Row 4
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.22"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -0.34%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  -0.03%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4898"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  +1.98%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2972"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = "  +1.76%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06872"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  +1.20%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.19"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  -0.54%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "108.01"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  -3.29%  "

# Row 12
$ws.Range("D12").Value = "1.957.51"
$ws.Range("E12").Value = "  +0.70%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.494"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  +0.25%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.7107"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  +3.84%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "283.27"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  -3.14%  "

# Row 17
$ws.Range("D17").Value = "31.104.69"
$ws.Range("E17").Value = "  +1.48%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.31"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  +0.78%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007785"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  +1.61%  "

# Row 20
$ws.Range("D20").Value = "2.209.13"
$ws.Range("E20").Value = "  +0.80%  "

# Row 21
$ws.Range("E21").Value = "  -0.02%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.535"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  -1.85%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  -0.01%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.522"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  -0.60%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.854"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  +1.14%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "169.77"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  +0.95%  "

# Row 27
$ws.Range("E27").Value = "  -0.59%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.231"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  +2.57%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1059"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  -1.76%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.429"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  -0.24%  "

# Row 31
$ws.Range("E31").Value = "  -0.32%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.596"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  -2.00%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.513"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  +0.83%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04996"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  -0.90%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7608"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  -0.86%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.186"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  +2.68%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.731"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  -0.13%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02034"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  -1.49%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.708"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  +0.38%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.171"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  +6.31%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.529"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  +10.04%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "74.71"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  +8.23%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4529"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  +1.89%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8887"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  +2.24%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "109.59"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  -0.91%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.126"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  +10.83%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.001"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  -0.22%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "985.97"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  +9.74%  "

# Row 49
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1271"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  +1.90%  "

# Row 50
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.413"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  +1.23%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.2585"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  +3.19%  "
